$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Fix the "Unit" text in G13 - remove stray spaces and trailing space
$ws.Range("G13").Value = "t/cap/year"

# Update the "Zeile im Spreadsheet" (row in spreadsheet) column K values
$ws.Range("K8").Value = 547
$ws.Range("K9").Value = 415
$ws.Range("K13").Value = 342
$ws.Range("K15").Value = 569
$ws.Range("K20").Value = 337
$ws.Range("K22").Value = 567
$ws.Range("K30").Value = 549
$ws.Range("K31").Value = 414
$ws.Range("K35").Value = 561
$ws.Range("K41").Value = 556
$ws.Range("K44").Value = 372
$ws.Range("K47").Value = 558
$ws.Range("K50").Value = 374
$ws.Range("K54").Value = 554
$ws.Range("K58").Value = 366
$ws.Range("K63").Value = 552
$ws.Range("K65").Value = 563
$ws.Range("K68").Value = 355
$ws.Range("K69").Value = 322
$ws.Range("K70").Value = 323
$ws.Range("K71").Value = 565
